# Update the "predictions" column (E2:E11) on Sheet1: the series was
# shifted down by one row (a new leading value was inserted and the old
# trailing value dropped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(30200, 30367, 30567, 30978, 30920, 31219, 31469, 31499, 31711, 31994)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Range("E$row").Value = $values[$i]
}

# Match the saved selection/active cell recorded in the workbook.
$ws.Range("E12").Select()
